# Vinay - ULP Changes - 0628
# Add a "State" column to each state worksheet, populated with that
# worksheet's own state abbreviation (matches the sheet/tab name).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $state = $ws.Name

    $ws.Range("C1").Value = "State"
    $ws.Range("C2").Value = $state

    $ws.Range("D2").Select()
}
